$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 130826010
$ws.Range("B2").Value = 91808
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = "Ullticka"
$ws.Range("G2").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H2").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M2").ClearContents()
$ws.Range("P2").Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Range("Q2").Value = 460971
$ws.Range("R2").Value = 7039688
$ws.Range("Z2").Value = "10:47"
$ws.Range("AB2").Value = "10:47"
$ws.Range("AC2").ClearContents()

# Row 3
$ws.Range("A3").Value = 130825823
$ws.Range("B3").Value = 57881
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = "Spillkråka"
$ws.Range("G3").Value = "Dryocopus martius"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("M3").Value = "äldre spår"
$ws.Range("P3").Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Range("Q3").Value = 460947
$ws.Range("R3").Value = 7039711
$ws.Range("Z3").Value = "10:38"
$ws.Range("AB3").Value = "10:38"
$ws.Range("AC3").Value = "Födosökshål på äldre döende gran."

# Row 4
$ws.Range("A4").Value = 130826784
$ws.Range("P4").Value = "Brännan, Kälom, Offerdal, Jmt"
$ws.Range("Q4").Value = 461233
$ws.Range("R4").Value = 7039438
$ws.Range("S4").Value = 10
$ws.Range("Z4").Value = "11:37"
$ws.Range("AB4").Value = "11:37"
$ws.Range("AC4").Value = "Födosök barkfläk"

# Row 5
$ws.Range("A5").Value = 130825852
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("M5").Value = "färska spår"
$ws.Range("Q5").Value = 460952
$ws.Range("R5").Value = 7039723
$ws.Range("S5").Value = 15
$ws.Range("Z5").Value = "10:42"
$ws.Range("AB5").Value = "10:42"
$ws.Range("AC5").Value = "Barkfläkta grövre och klenare granar."

# Row 22
$ws.Range("A22").Value = 130826355
$ws.Range("B22").Value = 92535
$ws.Range("D22").Value = "VU"
$ws.Range("E22").Value = 67
$ws.Range("F22").Value = "Sprickporing"
$ws.Range("G22").Value = "Diplomitoporus crustulinus"
$ws.Range("H22").Value = "(Bres.) Domański"
$ws.Range("P22").Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Range("Q22").Value = 461117
$ws.Range("R22").Value = 7039629
$ws.Range("S22").Value = 10
$ws.Range("Z22").Value = "11:10"
$ws.Range("AB22").Value = "11:10"
$ws.Range("AC22").Value = "På undersidan av lutande död gran."

# Row 23
$ws.Range("A23").Value = 130826438
$ws.Range("B23").Value = 79243
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("P23").Value = "Brännan, Brännan, Jmt"
$ws.Range("Q23").Value = 461220
$ws.Range("R23").Value = 7039590
$ws.Range("S23").Value = 25
$ws.Range("Z23").Value = "11:16"
$ws.Range("AB23").Value = "11:16"
$ws.Range("AC23").Value = "Rikligt i området"
